$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write cell values in the same order the original author entered them so that
# the shared-string table ends up populated in the same sequence as the target.
$ws.Range("B4").Value = "Examples with label 5 are embedded using TabNet based on Soil_Type, closely to each other (the cluster has a significant distribution value of label 5). However, VIME embeddings of these examples are mostly based on Wilderness_Area which splits the labels worse."
$ws.Range("B5").Value = "Examples with label 6 are closely embedded in TabNet based on Elevation mostly, unlike VIME, where they are scattered across all embedding space"
$ws.Range("B7").Value = "TabNet captures more complex dividing of the dataset, based on different columns, rather than VIME that splits mostly on a single column"
$ws.Range("D4").Value = "TabNet captures more complex dividing of the dataset, based on different columns, rather than VIME that splits based on less columns, therefore its patterns are simpler"
$ws.Range("B6").Value = "More insights in the manner of the above"
$ws.Range("D5").Value = "Trance examples are embedded using TabNet in separate cluster, whereas using VIME they are embedded in a cluster mixing multiple genres"
$ws.Range("D6").Value = "On examples where TabNet outperforms VIME with label Psytrance, it can be seen that VIME confuses psytrance with trance based on their instrumentalness column )have the same range and therefore embeds them closely), where TabNet splits them to different clusters"

# Columns auto-sized themselves (bestFit) in Excel after the new, longer
# strings were entered; reproduce the resulting widths as closely as possible.
$ws.Columns.Item(2).ColumnWidth = 236.77734375
$ws.Columns.Item(4).ColumnWidth = 150.44140625

$ws.Range("D7").Select()
$excel.ActiveWindow.ScrollColumn = 3
